$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 (new columns). Copy format from H1 (bold, border, centered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I and J, rows 2-37.
$data = @(
  @(6,7),
  @(8,8),
  @(6,6),
  @(5,6),
  @(9,9),
  @(8,8),
  @(9,9),
  @(8,8),
  @(7,7),
  @(9,9),
  @(6,6),
  @(9,9),
  @(9,9),
  @(9,9),
  @(6,6),
  @(6,6),
  @(7,7),
  @(6,6),
  @(3,4),
  @(5,5),
  @(8,8),
  @(6,6),
  @(6,6),
  @(5,5),
  @(6,6),
  @(7,7),
  @(7,7),
  @(8,8),
  @(3,4),
  @(3,3),
  @(7,7),
  @(7,7),
  @(6,6),
  @(6,6),
  @(6,6),
  @(1,1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 9).Value = $data[$i][0]
  $ws.Cells.Item($r, 10).Value = $data[$i][1]
}

Write-Host "Done"
